$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.792.50'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '2.319.10'
$ws.Range('E3').Value = '  +3.58%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '97.78'
$ws.Range('E5').Value = '  +6.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '272.69'
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.626'
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.45'
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.05'
$ws.Range('E12').Value = '  -2.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.106'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').Value = '2.658.45'
$ws.Range('E14').Value = '  +3.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.49'
$ws.Range('E15').Value = '  +2.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.875'
$ws.Range('E16').Value = '  +8.68%  '
$ws.Range('D17').Value = '2.330.20'
$ws.Range('E17').Value = '  +4.47%  '
$ws.Range('D18').Value = '43.750.76'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000109'
$ws.Range('E19').Value = '  +3.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.29'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '240.24'
$ws.Range('E22').Value = '  +2.69%  '
$ws.Range('E23').Value = '  -3.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.46'
$ws.Range('E24').Value = '  +4.77%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.55'
$ws.Range('E26').Value = '  +1.96%  '
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.51'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.31'
$ws.Range('E30').Value = '  -7.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.41'
$ws.Range('E31').Value = '  +7.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '175.06'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.48'
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('E35').Value = '  +2.85%  '
$ws.Range('E36').Value = '  -3.55%  '
$ws.Range('E37').Value = '  +3.05%  '
$ws.Range('E38').Value = '  +2.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.40'
$ws.Range('E39').Value = '  -3.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.243'
$ws.Range('E40').Value = '  +5.74%  '
$ws.Range('E41').Value = '  +9.08%  '
$ws.Range('E42').Value = '  +21.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.22'
$ws.Range('E43').Value = '  -4.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.96'
$ws.Range('E44').Value = '  -0.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.23'
$ws.Range('E45').Value = '  +10.27%  '
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E47').Value = '  +3.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '100.46'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.193'
$ws.Range('E50').Value = '  +16.89%  '
$ws.Range('D51').Value = '2.545.18'
$ws.Range('E51').Value = '  +3.50%  '
